$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 5140.3335
$ws.Cells.Item(62, 9).Value = 4658.6665
$ws.Cells.Item(62, 11).Value = 4658.6665
$ws.Cells.Item(62, 13).Value = -4034.6665

$ws.Cells.Item(65, 8).Value = 5140.3335
$ws.Cells.Item(65, 9).Value = 4658.6665
$ws.Cells.Item(65, 11).Value = 23293.3325
$ws.Cells.Item(65, 13).Value = -20173.3325

$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).ClearContents()

$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).ClearContents()

$ws.Cells.Item(116, 8).Value = 13529
$ws.Cells.Item(116, 9).Value = 24469.8
$ws.Cells.Item(116, 10).Value = 5714.143
$ws.Cells.Item(116, 11).Value = 24469.8
$ws.Cells.Item(116, 12).Value = 5714.143
$ws.Cells.Item(116, 13).Value = -21027.8
$ws.Cells.Item(116, 14).Value = -12598.143

$ws.Cells.Item(137, 8).Value = 9048.464
$ws.Cells.Item(137, 10).Value = 23873.715
$ws.Cells.Item(137, 12).Value = 71621.145
$ws.Cells.Item(137, 14).Value = -76721.145

$ws.Cells.Item(138, 8).Value = 1773669.1
$ws.Cells.Item(138, 9).Value = 1909.8
$ws.Cells.Item(138, 11).Value = 5729.4
$ws.Cells.Item(138, 13).Value = -589.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4833.625
$ws.Cells.Item(45, 9).Value = 1307
$ws.Cells.Item(45, 10).Value = 6949.6
$ws.Cells.Item(45, 11).Value = 1307
$ws.Cells.Item(45, 12).Value = 6949.6
$ws.Cells.Item(45, 13).Value = -930
$ws.Cells.Item(45, 14).Value = -7703.6

$ws.Cells.Item(46, 8).Value = 6257.143
$ws.Cells.Item(46, 9).Value = 6199.5
$ws.Cells.Item(46, 10).Value = 6334
$ws.Cells.Item(46, 11).Value = 6199.5
$ws.Cells.Item(46, 12).Value = 6334
$ws.Cells.Item(46, 13).Value = -5880.5
$ws.Cells.Item(46, 14).Value = -6972

$ws.Cells.Item(61, 8).Value = 6102.268
$ws.Cells.Item(61, 10).Value = 12927.4
$ws.Cells.Item(61, 12).Value = 12927.4
$ws.Cells.Item(61, 14).Value = -13351.4

$ws.Cells.Item(74, 8).Value = 3467.0408
$ws.Cells.Item(74, 9).Value = 1220.3235
$ws.Cells.Item(74, 11).Value = 1220.3235
$ws.Cells.Item(74, 13).Value = -346.3235

$ws.Cells.Item(77, 8).Value = 3467.0408
$ws.Cells.Item(77, 9).Value = 1220.3235
$ws.Cells.Item(77, 11).Value = 6101.6175
$ws.Cells.Item(77, 13).Value = -1733.6175

$ws.Cells.Item(97, 8).Value = 417.64
$ws.Cells.Item(97, 9).Value = 310.04166
$ws.Cells.Item(97, 10).Value = 3000
$ws.Cells.Item(97, 11).Value = 310.04166
$ws.Cells.Item(97, 12).Value = 3000
$ws.Cells.Item(97, 13).Value = 185.95834
$ws.Cells.Item(97, 14).Value = -3992

$ws.Cells.Item(102, 8).Value = 1633.7142
$ws.Cells.Item(102, 9).Value = 1581.1111
$ws.Cells.Item(102, 11).Value = 1581.1111
$ws.Cells.Item(102, 13).Value = 40.88889999999992

$ws.Cells.Item(122, 8).Value = 3401.625
$ws.Cells.Item(122, 9).Value = 2655
$ws.Cells.Item(122, 11).Value = 7965
$ws.Cells.Item(122, 13).Value = -5515

$ws.Cells.Item(136, 8).Value = 6102.268
$ws.Cells.Item(136, 10).Value = 12927.4
$ws.Cells.Item(136, 12).Value = 38782.2
$ws.Cells.Item(136, 14).Value = -43882.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 604.7646999999999
$ws.Cells.Item(80, 9).Value = 1021.2857
$ws.Cells.Item(80, 11).Value = 1021.2857
$ws.Cells.Item(80, 13).Value = -23.28570000000002

$ws.Cells.Item(83, 8).Value = 604.7646999999999
$ws.Cells.Item(83, 9).Value = 1021.2857
$ws.Cells.Item(83, 11).Value = 5106.4285
$ws.Cells.Item(83, 13).Value = -114.4285

$ws.Cells.Item(99, 8).Value = 4501.8887
$ws.Cells.Item(99, 9).Value = 3788.1428
$ws.Cells.Item(99, 11).Value = 3788.1428
$ws.Cells.Item(99, 13).Value = -2290.1428

$ws.Cells.Item(105, 8).Value = 2820
$ws.Cells.Item(105, 9).Value = 2194.25
$ws.Cells.Item(105, 11).Value = 2194.25
$ws.Cells.Item(105, 13).Value = -447.25

$ws.Cells.Item(134, 8).Value = 2794.432
$ws.Cells.Item(134, 9).Value = 2586.45
$ws.Cells.Item(134, 11).Value = 7759.349999999999
$ws.Cells.Item(134, 13).Value = -5224.349999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 157.30435
$ws.Cells.Item(7, 9).Value = 128.8421
$ws.Cells.Item(7, 10).Value = 292.5
$ws.Cells.Item(7, 11).Value = 128.8421
$ws.Cells.Item(7, 12).Value = 292.5
$ws.Cells.Item(7, 13).Value = -15.84209999999999
$ws.Cells.Item(7, 14).Value = -518.5

$ws.Cells.Item(99, 8).Value = 4409
$ws.Cells.Item(99, 9).Value = 4027
$ws.Cells.Item(99, 11).Value = 4027
$ws.Cells.Item(99, 13).Value = -2529

$ws.Cells.Item(126, 8).Value = 4409
$ws.Cells.Item(126, 9).Value = 4027
$ws.Cells.Item(126, 11).Value = 12081
$ws.Cells.Item(126, 13).Value = -9611

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 13681.818
$ws.Cells.Item(87, 10).Value = 13750
$ws.Cells.Item(87, 12).Value = 41250
$ws.Cells.Item(87, 14).Value = -43746

$ws.Cells.Item(90, 8).Value = 13681.818
$ws.Cells.Item(90, 10).Value = 13750
$ws.Cells.Item(90, 12).Value = 123750
$ws.Cells.Item(90, 14).Value = -136230

$ws.Cells.Item(132, 8).Value = 1359.1666
$ws.Cells.Item(132, 9).Value = 1114.6666
$ws.Cells.Item(132, 11).Value = 10031.9994
$ws.Cells.Item(132, 13).Value = -7501.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1840.2
$ws.Cells.Item(102, 9).Value = 1584.9231
$ws.Cells.Item(102, 11).Value = 1584.9231
$ws.Cells.Item(102, 13).Value = 37.07690000000002

$ws.Cells.Item(141, 8).Value = 66140
$ws.Cells.Item(141, 10).Value = 66140
$ws.Cells.Item(141, 12).Value = 66140
$ws.Cells.Item(141, 14).Value = -76500

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1820.9375
$ws.Cells.Item(22, 10).Value = 1945.5
$ws.Cells.Item(22, 12).Value = 1945.5
$ws.Cells.Item(22, 14).Value = -2535.5

$ws.Cells.Item(27, 8).Value = 1820.9375
$ws.Cells.Item(27, 10).Value = 1945.5
$ws.Cells.Item(27, 12).Value = 1945.5
$ws.Cells.Item(27, 14).Value = -2159.5

$ws.Cells.Item(40, 8).Value = 3422
$ws.Cells.Item(40, 9).Value = 3308.5264
$ws.Cells.Item(40, 10).Value = 4500
$ws.Cells.Item(40, 11).Value = 3308.5264
$ws.Cells.Item(40, 12).Value = 4500
$ws.Cells.Item(40, 13).Value = -3172.5264
$ws.Cells.Item(40, 14).Value = -4772

$ws.Cells.Item(122, 8).Value = 3178.4
$ws.Cells.Item(122, 10).Value = 3130.6667
$ws.Cells.Item(122, 12).Value = 9392.000100000001
$ws.Cells.Item(122, 14).Value = -14292.0001

$ws.Cells.Item(132, 8).Value = 4526.543
$ws.Cells.Item(132, 9).Value = 4433.35
$ws.Cells.Item(132, 10).Value = 4650.8
$ws.Cells.Item(132, 11).Value = 13300.05
$ws.Cells.Item(132, 12).Value = 13952.4
$ws.Cells.Item(132, 13).Value = -10770.05
$ws.Cells.Item(132, 14).Value = -19012.4

$ws.Cells.Item(136, 8).Value = 5312.1577
$ws.Cells.Item(136, 9).Value = 3963.2856
$ws.Cells.Item(136, 10).Value = 6099
$ws.Cells.Item(136, 11).Value = 11889.8568
$ws.Cells.Item(136, 12).Value = 18297
$ws.Cells.Item(136, 13).Value = -9339.856800000001
$ws.Cells.Item(136, 14).Value = -23397

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4959.7036
$ws.Cells.Item(81, 9).Value = 4900.8184
$ws.Cells.Item(81, 10).Value = 5218.8
$ws.Cells.Item(81, 11).Value = 9801.6368
$ws.Cells.Item(81, 12).Value = 10437.6
$ws.Cells.Item(81, 13).Value = -8740.6368
$ws.Cells.Item(81, 14).Value = -12559.6

$ws.Cells.Item(84, 8).Value = 4959.7036
$ws.Cells.Item(84, 9).Value = 4900.8184
$ws.Cells.Item(84, 10).Value = 5218.8
$ws.Cells.Item(84, 11).Value = 49008.184
$ws.Cells.Item(84, 12).Value = 52188
$ws.Cells.Item(84, 13).Value = -43704.184
$ws.Cells.Item(84, 14).Value = -62796

$ws.Cells.Item(96, 8).Value = 333334400
$ws.Cells.Item(96, 9).Value = 1599.5
$ws.Cells.Item(96, 10).Value = 1000000000
$ws.Cells.Item(96, 11).Value = 1599.5
$ws.Cells.Item(96, 12).Value = 1000000000
$ws.Cells.Item(96, 13).Value = -226.5
$ws.Cells.Item(96, 14).Value = -1000002746

$ws.Cells.Item(126, 8).Value = 2808.3635
$ws.Cells.Item(126, 10).Value = 3786.125
$ws.Cells.Item(126, 12).Value = 11358.375
$ws.Cells.Item(126, 14).Value = -16298.375

$ws.Cells.Item(132, 8).Value = 25073382
$ws.Cells.Item(132, 9).Value = 35816972
$ws.Cells.Item(132, 11).Value = 107450916
$ws.Cells.Item(132, 13).Value = -107448386
